$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 379, shifting all existing rows 379-430 down to 380-431.
$ws.Rows.Item(379).Insert()

# Populate the newly inserted row 379 with the new weekly record.
$ws.Cells.Item(379, 1).Value = 9
$ws.Cells.Item(379, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(379, 3).Value = "Metropolitana"
$ws.Cells.Item(379, 4).Value = 45131
$ws.Cells.Item(379, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(379, 5).Value = 13
$ws.Cells.Item(379, 6).Value = 100112021
$ws.Cells.Item(379, 7).Value = "Ají"
$ws.Cells.Item(379, 8).Value = "Inferno"
$ws.Cells.Item(379, 9).Value = "Primera"
$ws.Cells.Item(379, 10).Value = 70
$ws.Cells.Item(379, 11).Value = 13000
$ws.Cells.Item(379, 12).Value = 14000
$ws.Cells.Item(379, 13).Value = 13500
$ws.Cells.Item(379, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(379, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(379, 16).Value = 1350
$ws.Cells.Item(379, 17).Value = 10
$ws.Cells.Item(379, 18).Value = "Hortaliza"
